$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for new column F (match style/formatting of the other header cells)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$timestamps = @(
    "2021-10-05 10:51:06.737669",
    "2021-10-05 10:51:06.737679",
    "2021-10-05 10:51:06.737683",
    "2021-10-05 10:51:06.737685",
    "2021-10-05 10:51:06.737688",
    "2021-10-05 10:51:06.737691",
    "2021-10-05 10:51:06.737693",
    "2021-10-05 10:51:06.737696",
    "2021-10-05 10:51:06.737699",
    "2021-10-05 10:51:06.737701",
    "2021-10-05 10:51:06.737704",
    "2021-10-05 10:51:06.737706",
    "2021-10-05 10:51:06.737709",
    "2021-10-05 10:51:06.737711",
    "2021-10-05 10:51:06.737714",
    "2021-10-05 10:51:06.737716",
    "2021-10-05 10:51:06.737719",
    "2021-10-05 10:51:06.737722",
    "2021-10-05 10:51:06.737724",
    "2021-10-05 10:51:06.737727",
    "2021-10-05 10:51:06.737729",
    "2021-10-05 10:51:06.737732",
    "2021-10-05 10:51:06.737734",
    "2021-10-05 10:51:06.737737",
    "2021-10-05 10:51:06.737740",
    "2021-10-05 10:51:06.737743",
    "2021-10-05 10:51:06.737745",
    "2021-10-05 10:51:06.737748",
    "2021-10-05 10:51:06.737751",
    "2021-10-05 10:51:06.737753",
    "2021-10-05 10:51:06.737756",
    "2021-10-05 10:51:06.737758",
    "2021-10-05 10:51:06.737761",
    "2021-10-05 10:51:06.737764",
    "2021-10-05 10:51:06.737766",
    "2021-10-05 10:51:06.737769"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
